$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 27, shifting rows 27:150 down to 28:151
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record
$ws.Range("A27").Value = 9
$ws.Range("B27").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 45250
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 100114007
$ws.Range("G27").Value = "Jengibre"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 700
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 20000
$ws.Range("M27").Value = 19000
$ws.Range("N27").Value = "$/caja 13 kilos"
$ws.Range("O27").Value = "Perú"
$ws.Range("P27").Value = 1462
$ws.Range("Q27").Value = 13
$ws.Range("R27").Value = "Hortaliza"
